# refactor: update 2a, 2b, 3 to collect prior results
#
# Rows 8-13 previously stored the PCR/Seq "prior results" dates (K, L) as
# literal text strings ("04/22/2023" / "05/19/2023"). Switch them to real
# date values so they collect/compute like the other rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K8:K13 -> 04/22/2023 as a real date serial (was text "04/22/2023")
$ws.Range("K8:K13").Value = 45038

# L8:L13 -> 05/19/2023 as a real date serial (was text "05/19/2023")
$ws.Range("L8:L13").Value = 45065

# Selection moves to the newly-updated prior-results range.
$ws.Range("L8:L13").Select()
